# Update BOM price/subtotal values to latest supplier quotes.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4
$ws.Range("H4").Value = 0.01

# Row 6
$ws.Range("G6").Value = 0.0257
$ws.Range("H6").Value = 0.257

# Row 8
$ws.Range("G8").Value = 0.0483
$ws.Range("H8").Value = 0.0483

# Row 9
$ws.Range("G9").Value = 0.021
$ws.Range("H9").Value = 0.042

# Row 10
$ws.Range("G10").Value = 0.098
$ws.Range("H10").Value = 0.098

# Row 11
$ws.Range("G11").Value = 0.1059
$ws.Range("H11").Value = 0.2118

# Row 14
$ws.Range("G14").Value = 0.011
$ws.Range("H14").Value = 0.11

# Row 21
$ws.Range("G21").Value = 2.5
$ws.Range("H21").Value = 2.5

# Row 22
$ws.Range("G22").Value = 0.015
$ws.Range("H22").Value = 0.15

# Row 23
$ws.Range("G23").Value = 0.013
$ws.Range("H23").Value = 0.13

# Row 25
$ws.Range("G25").Value = 0.012
$ws.Range("H25").Value = 0.12

# Row 26
$ws.Range("G26").Value = 0.024
$ws.Range("H26").Value = 0.024

# Row 27
$ws.Range("G27").Value = 0.013
$ws.Range("H27").Value = 0.052

# Row 29
$ws.Range("G29").Value = 0.01
$ws.Range("H29").Value = 0.1

# Row 33
$ws.Range("G33").Value = 4.85
$ws.Range("H33").Value = 4.85

# Row 36
$ws.Range("G36").Value = 1.71
$ws.Range("H36").Value = 1.71

# Row 39 (new values added where there were none before)
$ws.Range("G39").Value = 0.41
$ws.Range("H39").Value = 0.41

# Row 40
$ws.Range("G40").Value = 0.9533
$ws.Range("H40").Value = 0.9533
